# Add script for "add loan type" screen to the Institute sheet (TestData.xlsx)
#
# This appends five new columns (BH:BL) to the "Institute" worksheet:
#   BH = HolidayType
#   BI = LoanType
#   BJ = DraftNeeded
#   BK = LoanTypeCreditLimit
#   BL = LoanTypeCashLimit
# with header labels on row 1 and sample data on rows 2-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Institute")

# --- Header row (row 1) ---------------------------------------------------
$headers = @("HolidayType", "LoanType", "DraftNeeded", "LoanTypeCreditLimit", "LoanTypeCashLimit")
$headerCols = @("BH", "BI", "BJ", "BK", "BL")

for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $headers[$i]
}

# Copy the existing header formatting (fill/font/border) onto the new header cells
$ws.Range("BG1").Copy() | Out-Null
$ws.Range("BH1:BL1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows (rows 2-11) -------------------------------------------------
# BI3 (row 3) carries a different value from the rest of the rows
$loanTypeByRow = @{
    2  = "Loan below credit limit [LOANCR]"
    3  = "Retail Transaction to Loan [LOANPUR]"
    4  = "Loan below credit limit [LOANCR]"
    5  = "Loan below credit limit [LOANCR]"
    6  = "Loan below credit limit [LOANCR]"
    7  = "Loan below credit limit [LOANCR]"
    8  = "Loan below credit limit [LOANCR]"
    9  = "Loan below credit limit [LOANCR]"
    10 = "Loan below credit limit [LOANCR]"
    11 = "Loan below credit limit [LOANCR]"
}

for ($row = 2; $row -le 11; $row++) {
    $ws.Range("BH$row").Value = "Holiday [H]"
    $ws.Range("BI$row").Value = $loanTypeByRow[$row]
    $ws.Range("BJ$row").Value = "Check"
    $ws.Range("BK$row").Value = "Check"
    $ws.Range("BL$row").Value = "Check"

    # Copy the existing data-row formatting onto the new cells
    $ws.Range("BG$row").Copy() | Out-Null
    $ws.Range("BH" + $row + ":BL" + $row).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# --- Column width for the new last column (BL / column 64) -----------------
$ws.Columns.Item(64).ColumnWidth = 16.7

# --- Selection / active cell ------------------------------------------------
$ws.Activate()
$ws.Range("BK3").Select() | Out-Null

# --- Workbook calculation option (best effort) ------------------------------
$excel.MultiThreadedCalculation.Enabled = $false
